$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value to a cell while preventing Excel from
# auto-converting numeric-looking text (e.g. "1.00") into a real number.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "64.982.13"
$ws.Range("E2").Value = "  -6.80%  "
# Row 3
$ws.Range("D3").Value = "3.274.30"
$ws.Range("E3").Value = "  -7.76%  "
# Row 4
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.16%  "
# Row 5
Set-TextValue $ws.Range("D5") "550.12"
$ws.Range("E5").Value = "  -7.06%  "
# Row 6
Set-TextValue $ws.Range("D6") "177.45"
$ws.Range("E6").Value = "  -9.55%  "
# Row 7
$ws.Range("E7").Value = "  +0.07%  "
# Row 8
Set-TextValue $ws.Range("D8") "0.584"
$ws.Range("E8").Value = "  -4.44%  "
# Row 9
$ws.Range("D9").Value = "3.279.79"
$ws.Range("E9").Value = "  -7.26%  "
# Row 10
Set-TextValue $ws.Range("D10") "0.182"
$ws.Range("E10").Value = "  -13.07%  "
# Row 11
Set-TextValue $ws.Range("D11") "0.578"
$ws.Range("E11").Value = "  -7.53%  "
# Row 12
Set-TextValue $ws.Range("D12") "46.57"
$ws.Range("E12").Value = "  -11.66%  "
# Row 13
Set-TextValue $ws.Range("D13") "0.0000260"
$ws.Range("E13").Value = "  -9.98%  "
# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.819.59"
$ws.Range("E14").Value = "  -7.26%  "
# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D15") "8.43"
$ws.Range("E15").Value = "  -8.75%  "
# Row 16
Set-TextValue $ws.Range("D16") "598.83"
$ws.Range("E16").Value = "  -9.07%  "
# Row 17
Set-TextValue $ws.Range("D17") "17.81"
$ws.Range("E17").Value = "  -3.32%  "
# Row 18
$ws.Range("D18").Value = "65.081.77"
$ws.Range("E18").Value = "  -6.50%  "
# Row 19
$ws.Range("E19").Value = "  -4.34%  "
# Row 20
$ws.Range("D20").Value = "3.292.32"
$ws.Range("E20").Value = "  -7.50%  "
# Row 21
Set-TextValue $ws.Range("D21") "11.25"
$ws.Range("E21").Value = "  -10.61%  "
# Row 22
Set-TextValue $ws.Range("D22") "0.891"
$ws.Range("E22").Value = "  -7.32%  "
# Row 23
Set-TextValue $ws.Range("D23") "17.22"
$ws.Range("E23").Value = "  -5.26%  "
# Row 24
Set-TextValue $ws.Range("D24") "101.00"
$ws.Range("E24").Value = "  -3.43%  "
# Row 25
Set-TextValue $ws.Range("D25") "4.91"
$ws.Range("E25").Value = "  -8.65%  "
# Row 26
Set-TextValue $ws.Range("D26") "3.93"
$ws.Range("E26").Value = "  -10.65%  "
# Row 27
Set-TextValue $ws.Range("D27") "5.96"
$ws.Range("E27").Value = "  -1.30%  "
# Row 28
Set-TextValue $ws.Range("D28") "2.64"
$ws.Range("E28").Value = "  -9.51%  "
# Row 29
Set-TextValue $ws.Range("D29") "9.21"
$ws.Range("E29").Value = "  -9.70%  "
# Row 30
Set-TextValue $ws.Range("D30") "8.53"
$ws.Range("E30").Value = "  -11.08%  "
# Row 31
Set-TextValue $ws.Range("D31") "29.98"
$ws.Range("E31").Value = "  -9.97%  "
# Row 32
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D32") "6.14"
$ws.Range("E32").Value = "  -9.63%  "
# Row 33
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D33") "3.69"
$ws.Range("E33").Value = "  -15.94%  "
# Row 34
Set-TextValue $ws.Range("D34") "10.88"
$ws.Range("E34").Value = "  -7.54%  "
# Row 35
$ws.Range("D35").Value = "3.780.75"
$ws.Range("E35").Value = "  +0.77%  "
# Row 36
Set-TextValue $ws.Range("D36") "0.103"
$ws.Range("E36").Value = "  -7.11%  "
# Row 37
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D37") "0.999"
$ws.Range("E37").Value = "  +0.03%  "
# Row 38
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D38") "524.06"
$ws.Range("E38").Value = "  +2.80%  "
# Row 39
Set-TextValue $ws.Range("D39") "55.54"
$ws.Range("E39").Value = "  -10.16%  "
# Row 40
Set-TextValue $ws.Range("D40") "3.37"
$ws.Range("E40").Value = "  -10.17%  "
# Row 41
$ws.Range("D41").Value = "0.0₃0701"
$ws.Range("E41").Value = "  -13.92%  "
# Row 42
Set-TextValue $ws.Range("D42") "2.61"
$ws.Range("E42").Value = "  -10.73%  "
# Row 43
Set-TextValue $ws.Range("D43") "0.124"
$ws.Range("E43").Value = "  -7.87%  "
# Row 44
Set-TextValue $ws.Range("D44") "0.333"
$ws.Range("E44").Value = "  -10.23%  "
# Row 45
Set-TextValue $ws.Range("D45") "31.35"
$ws.Range("E45").Value = "  -9.95%  "
# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D46") "3.22"
$ws.Range("E46").Value = "  -5.86%  "
# Row 47
$ws.Range("B47").Value = "CoreDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-TextValue $ws.Range("D47") "3.11"
$ws.Range("E47").Value = "  +14.44%  "
# Row 48
Set-TextValue $ws.Range("D48") "0.0404"
$ws.Range("E48").Value = "  -10.88%  "
# Row 49
Set-TextValue $ws.Range("D49") "0.128"
$ws.Range("E49").Value = "  -6.33%  "
# Row 50
Set-TextValue $ws.Range("D50") "2.56"
$ws.Range("E50").Value = "  -10.92%  "
# Row 51
Set-TextValue $ws.Range("D51") "1.00"
$ws.Range("E51").Value = "  +0.13%  "
